# Update "Lương" sheet with new values reflecting an extra day of salary advance detail
# Commit: Thêm chi tiết về ứng lương vào báo cảo tổng hợp

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lương")

$ws.Range("B12").Value = 24
$ws.Range("B13").Value = 3428571.428571429
$ws.Range("B32").Value = 1128571.428571429
$ws.Range("B34").Value = 1128571.428571429


